$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.872.84"
$ws.Range("E2").Value = "  +2.56%  "
$ws.Range("D3").Value = "2.227.42"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "270.37"
$ws.Range("E5").Value = "  +5.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.23"
$ws.Range("E6").Value = "  +11.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  +1.54%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.607"
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.37"
$ws.Range("E10").Value = "  +9.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0923"
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.54"
$ws.Range("E12").Value = "  +7.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.105"
$ws.Range("E13").Value = "  +1.86%  "
$ws.Range("D14").Value = "2.559.18"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.73"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").Value = "2.225.58"
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.792"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "43.857.32"
$ws.Range("E18").Value = "  +2.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000104"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.03"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.99"
$ws.Range("E21").Value = "  -1.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.39"
$ws.Range("E22").Value = "  +4.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.66"
$ws.Range("E23").Value = "  +1.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.92"
$ws.Range("E24").Value = "  -4.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.70"
$ws.Range("E25").Value = "  +22.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.83"
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("E28").Value = "  +5.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.63"
$ws.Range("E29").Value = "  -7.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.25"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.30"
$ws.Range("E31").Value = "  +1.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0902"
$ws.Range("E32").Value = "  +3.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.61"
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.43"
$ws.Range("E34").Value = "  +3.65%  "
$ws.Range("E35").Value = "  +2.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.111"
$ws.Range("E36").Value = "  +3.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0358"
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.38"
$ws.Range("E38").Value = "  +0.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.36"
$ws.Range("E39").Value = "  +19.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.47"
$ws.Range("E40").Value = "  -6.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "65.16"
$ws.Range("E41").Value = "  +8.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.11"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.204"
$ws.Range("E43").Value = "  +1.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.32"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0991"
$ws.Range("E45").Value = "  +1.50%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.96"
$ws.Range("E46").Value = "  -2.93%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.33"
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("E48").Value = "  +6.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.13"
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.444"
$ws.Range("E50").Value = "  -5.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.51"
$ws.Range("E51").Value = "  +4.61%  "
